$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 623 (the "猫" post row) -- this shifts all rows below it up by one.
$ws.Rows.Item(623).Delete()
